# repull data, push all data, mean calculation
# Update column F (dSF) values for the data rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -4
    4  = -4
    6  = 6
    8  = 1
    9  = 3
    10 = 5
    12 = -4
    13 = 8
    14 = -1
    15 = -2
    16 = -3
    17 = 4
    18 = 4
    20 = 9
    21 = -2
    22 = -3
    24 = 2
    25 = 1
    28 = 3
    29 = -3
    30 = 3
    31 = 3
    32 = 1
    33 = 5
    34 = 5
    35 = -2
    36 = 3
    37 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
